$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'262.14"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'22.92"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Value = "'6.195"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'0.06108"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'6.742"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'3.458"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'1.362"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.7986"
$ws.Range('D9').Style = 'Normal'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = "'0.1586"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = "'0.08078"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = "'0.03432"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03090"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09323"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').Value = "'3.847"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = "'0.001704"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = "'0.04837"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = "'0.0006137"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '17OneONEWorstin24h'
$ws.Range('D19').Value = "'0.006183"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'0.001093"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'0.003548"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Value = "'3.711"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'2.239"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Value = "'0.1251"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0003200"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D40').Value = "'0.04591"
$ws.Range('D40').Style = 'Normal'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = "'0.1119"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('D42').Value = "'0.003130"
$ws.Range('D42').Style = 'Normal'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = "'0.003331"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('D46').Value = "'0.00005921"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Value = "'0.6996"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'0.09295"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '48BOLOBOLO'
$ws.Range('D50').Value = "'0.00002099"
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = "'0.01009"
$ws.Range('D51').Style = 'Normal'
